$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (previously the last row) switches from the "date-only" style to
# the "date+time" style used by all the other data rows.
$ws.Range("A29").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 30 with the next day's data, using the "date-only" style that
# the previous last row (29) used to have.
$ws.Range("A30").Value = 45979
$ws.Range("B30").Value = 65
$ws.Range("C30").Value = 73
$ws.Range("D30").Value = 75

$ws.Range("A30").NumberFormat = "YYYY-MM-DD"
